# Add team record (Wins/Losses/Ties) columns to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, bordered, centered) from an existing header
# cell (AC1) onto the new header cells so they match the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows 2-44: every player/row gets the same team record.
$lastRow = 44
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 98  # AD = column 30 -> Wins
    $ws.Cells.Item($r, 31).Value = 64  # AE = column 31 -> Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF = column 32 -> Ties
}
